# Auto-generated edit script updating the cryptos worksheet values
# per the commit "Updated cryptos list ... with GitHub Actions".
# All touched cells hold text (coin names, links, formatted prices and
# percentages), so each is forced to Text format before the new value is
# written -- this stops Excel's COM layer from "helpfully" reinterpreting
# strings like "529.60" or "1.10" as numbers and dropping trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = 'D2'; Value = '57.750.20' }
    @{ Cell = 'E2'; Value = '  -4.70%  ' }
    @{ Cell = 'D3'; Value = '3.164.43' }
    @{ Cell = 'E3'; Value = '  -5.38%  ' }
    @{ Cell = 'E4'; Value = '  +0.07%  ' }
    @{ Cell = 'D5'; Value = '529.60' }
    @{ Cell = 'E5'; Value = '  -6.61%  ' }
    @{ Cell = 'D6'; Value = '134.84' }
    @{ Cell = 'E6'; Value = '  -8.03%  ' }
    @{ Cell = 'E7'; Value = '  -0.03%  ' }
    @{ Cell = 'D8'; Value = '3.164.94' }
    @{ Cell = 'E8'; Value = '  -5.36%  ' }
    @{ Cell = 'E9'; Value = '  -6.57%  ' }
    @{ Cell = 'D10'; Value = '7.19' }
    @{ Cell = 'E10'; Value = '  -8.98%  ' }
    @{ Cell = 'E11'; Value = '  -8.14%  ' }
    @{ Cell = 'E12'; Value = '  -4.70%  ' }
    @{ Cell = 'D13'; Value = '3.709.13' }
    @{ Cell = 'E13'; Value = '  -5.22%  ' }
    @{ Cell = 'E14'; Value = '  -1.53%  ' }
    @{ Cell = 'E15'; Value = '  -6.39%  ' }
    @{ Cell = 'D16'; Value = '3.166.34' }
    @{ Cell = 'E16'; Value = '  -5.04%  ' }
    @{ Cell = 'D17'; Value = '57.693.08' }
    @{ Cell = 'E17'; Value = '  -4.79%  ' }
    @{ Cell = 'E18'; Value = '  -8.59%  ' }
    @{ Cell = 'D19'; Value = '5.85' }
    @{ Cell = 'E19'; Value = '  -6.62%  ' }
    @{ Cell = 'D20'; Value = '13.20' }
    @{ Cell = 'E20'; Value = '  -9.22%  ' }
    @{ Cell = 'E21'; Value = '  -8.98%  ' }
    @{ Cell = 'D22'; Value = '349.72' }
    @{ Cell = 'E22'; Value = '  -7.09%  ' }
    @{ Cell = 'E23'; Value = '  +0.14%  ' }
    @{ Cell = 'D24'; Value = '69.77' }
    @{ Cell = 'E24'; Value = '  -6.64%  ' }
    @{ Cell = 'E25'; Value = '  -8.04%  ' }
    @{ Cell = 'D26'; Value = '3.294.13' }
    @{ Cell = 'E26'; Value = '  -5.61%  ' }
    @{ Cell = 'D27'; Value = '0.0₃0968' }
    @{ Cell = 'E27'; Value = '  -10.35%  ' }
    @{ Cell = 'E28'; Value = '  -3.89%  ' }
    @{ Cell = 'E29'; Value = '  -0.26%  ' }
    @{ Cell = 'E30'; Value = '  -5.25%  ' }
    @{ Cell = 'D31'; Value = '0.998' }
    @{ Cell = 'E31'; Value = '  -0.19%  ' }
    @{ Cell = 'E32'; Value = '  -9.08%  ' }
    @{ Cell = 'D33'; Value = '6.98' }
    @{ Cell = 'E33'; Value = '  -9.11%  ' }
    @{ Cell = 'D34'; Value = '21.70' }
    @{ Cell = 'E34'; Value = '  -5.19%  ' }
    @{ Cell = 'E35'; Value = '  -5.36%  ' }
    @{ Cell = 'D36'; Value = '4.99' }
    @{ Cell = 'E36'; Value = '  -5.66%  ' }
    @{ Cell = 'D37'; Value = '159.28' }
    @{ Cell = 'E37'; Value = '  -4.80%  ' }
    @{ Cell = 'E38'; Value = '  -7.89%  ' }
    @{ Cell = 'E39'; Value = '  -8.60%  ' }
    @{ Cell = 'D40'; Value = '26.33' }
    @{ Cell = 'E40'; Value = '  -5.95%  ' }
    @{ Cell = 'D42'; Value = '3.194.98' }
    @{ Cell = 'E42'; Value = '  -5.45%  ' }
    @{ Cell = 'D43'; Value = '40.36' }
    @{ Cell = 'E43'; Value = '  -4.35%  ' }
    @{ Cell = 'B44'; Value = 'ONDO' }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo' }
    @{ Cell = 'D44'; Value = '1.10' }
    @{ Cell = 'E44'; Value = '  -3.66%  ' }
    @{ Cell = 'B45'; Value = 'Mantle' }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt' }
    @{ Cell = 'D45'; Value = '0.699' }
    @{ Cell = 'E45'; Value = '  -7.72%  ' }
    @{ Cell = 'E46'; Value = '  -6.81%  ' }
    @{ Cell = 'E47'; Value = '  +0.03%  ' }
    @{ Cell = 'D48'; Value = '1.48' }
    @{ Cell = 'E48'; Value = '  -8.08%  ' }
    @{ Cell = 'D49'; Value = '2.275.82' }
    @{ Cell = 'E49'; Value = '  -7.26%  ' }
    @{ Cell = 'E50'; Value = '  -6.54%  ' }
    @{ Cell = 'D51'; Value = '20.74' }
    @{ Cell = 'E51'; Value = '  -7.27%  ' }
)

foreach ($change in $changes) {
    $range = $ws.Range($change.Cell)
    $range.NumberFormat = "@"
    $range.Value = $change.Value
}
